# Refresh the crypto price/volume table on Sheet1 (rows 2-51) with the
# latest scrape. Every value in columns D/E is stored as TEXT in this
# workbook (e.g. "70.035.61", "  -0.71%  "), so plain numeric-looking
# values are written with a leading apostrophe to stop Excel from
# auto-converting them to numbers, then NumberFormat is restored to
# "General" so the cell keeps its original (default) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.035.61"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.504.76"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'604.64"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'172.45"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "3.499.16"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").Value = "'7.25"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +7.34%  "
$ws.Range("D12").Value = "'0.586"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'46.01"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "4.069.82"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'8.35"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "'613.79"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "3.502.56"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "70.038.15"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "'17.52"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "'0.879"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'9.15"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -8.42%  "
$ws.Range("D24").Value = "'98.67"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "'15.52"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").Value = "'33.74"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'9.01"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").Value = "'8.05"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.30%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'632.66"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +11.40%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "'1.27"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").Value = "'0.0999"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").Value = "'10.74"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +4.38%  "
$ws.Range("D39").Value = "'3.47"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "'56.67"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").Value = "3.358.01"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "0.0₃0733"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").Value = "'0.311"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -5.35%  "
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").Value = "'31.86"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'133.33"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -0.02%  "
